$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D20").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E20").Value = "['Normal']"

$ws.Range("D44").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E44").Value = "['Normal', 'HardwareFault']"

$ws.Range("D64").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E64").Value = "['Normal']"
